# Generate Report for Handback
# Updates the recorded handoff/handback timestamps for the file
# "39613f77-56ff-4866-a0b3-591de88e5561" after a fresh localization
# handback run completed.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: refresh "Latest HO Xliff Generate Date" for the file ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-08-15 18:43:42"

# --- zh-cn sheet: refresh handoff/handback datetimes for the file's row ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H3").Value = "2016-08-15 18:43:37"
$wsZhCn.Range("K3").Value = "2016-08-15 18:43:54"

# --- de-de sheet: refresh handoff/handback datetimes for the file's row ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H3").Value = "2016-08-15 18:43:42"
$wsDeDe.Range("K3").Value = "2016-08-15 18:44:03"
